$d = $word.ActiveDocument

# NOTE: this runtime's Range.Find.Execute ignores the boundaries of the
# Range it's called on and searches/replaces across the *whole* document
# (confirmed experimentally - even Cell(r,c).Range.Find.Execute with
# Wrap=wdFindStop touches the first document-order match, not the one in
# that cell). Several of the new cell values below reproduce text that
# used to live in a *different* cell (e.g. "72÷9=" and "75÷8=" each
# appear twice), so a naive Find/Replace-based sequence quietly clobbers
# the wrong cell. Assigning directly to Range.Text is properly scoped to
# that Range, so it's used everywhere instead.

# Title line above the table.
$d.Paragraphs(1).Range.Text = "2025-01-20 Monday"

# Table cells, addressed by (row, column). Rows 2-4, 6-8, 10-12, 14-16,
# 18-20 are blank spacer rows; the five data rows are 1, 5, 9, 13, 17.
$tbl = $d.Tables(1)

$newGrid = @(
    @(1, @("43÷8=", "91÷4=", "51÷2=", "72÷9=", "93÷7=")),
    @(5, @("16÷7=", "12÷5=", "84÷5=", "75÷8=", "75÷8=")),
    @(9, @("29÷9=", "85÷4=", "95÷6=", "52÷8=", "29÷8=")),
    @(13, @("89÷5=", "76÷3=", "96÷9=", "31÷9=", "53÷5=")),
    @(17, @("87÷6=", "74÷4=", "60÷2=", "64÷8=", "54÷3="))
)

foreach ($entry in $newGrid) {
    $row = $entry[0]
    $vals = $entry[1]
    for ($col = 1; $col -le $vals.Count; $col++) {
        $cell = $tbl.Cell($row, $col)
        $cell.Range.Text = $vals[$col - 1]
    }
}
